$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "custom accuracy" - reduce the stored precision of row 5's measurements
#    (B5:AH5) from ~3 decimals down to 2 decimals, the same way Excel's
#    ROUND(x, 2) would (half-away-from-zero).
# ---------------------------------------------------------------------------
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = [Math]::Round([double]$v, 2)
    }
}

# ---------------------------------------------------------------------------
# 2) Drop the last data row (row 6) - part of the "1000 data points" re-export
#    that shortened this particular sheet by one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3) The narrower values mean several columns no longer need to be as wide;
#    set each column's width back to what autofit would have produced for
#    the new, shorter content (sheet column width is stored in "characters"
#    plus ~5/6 of a character of padding, so subtract that back out here).
# ---------------------------------------------------------------------------
$padding = 5.0 / 6.0
$targetWidths = @{
    2  = 7
    3  = 6
    5  = 7
    6  = 7
    8  = 7
    9  = 7
    12 = 7
    13 = 7
    15 = 7
    16 = 7
    20 = 7
    21 = 7
    23 = 7
    24 = 7
    26 = 7
    29 = 6
    30 = 7
    32 = 7
    34 = 7
}
foreach ($colIndex in $targetWidths.Keys) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetWidths[$colIndex] - $padding
}

Write-Output "done"
